$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update candidate personal info
$ws.Range("B2").Value = "ANG YI LING"
$ws.Range("C2").Value = "014-6263882"
$ws.Range("D2").Value = "yi_ling13@hotmail.com"

# "local" column changes from N/A to No; "expected_salary" stays N/A
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "N/A"

# Previous job roles
$ws.Range("G2").Value = "[{'job_title': 'SENIOR CORPORATE TRAVEL CONSULTANT', 'job_company': 'FCM TRAVEL SOLUTIONS MALAYSIA', 'Industries': 'Travel', 'start_date': '2018-11', 'end_date': '2021-07', 'job_location': 'Malaysia', 'job_duration': '2 years 9 months'}, {'job_title': 'SENIOR CORPORATE TRAVEL CONSULTANT', 'job_company': 'FCM TRAVEL SINGAPORE', 'Industries': 'Travel', 'start_date': '2021-08', 'end_date': '2021-11', 'job_location': 'Singapore', 'job_duration': '3 months'}, {'job_title': 'CORPORATE TRAVEL CONSULTANT', 'job_company': 'HOLIDAY TOURS SDN BHD', 'Industries': 'Travel', 'start_date': '2016-07', 'end_date': '2018-04', 'job_location': 'Kuala Lumpur', 'job_duration': '1 year 9 months'}]"

# Education background
$ws.Range("I2").Value = "[{'field_of_study': 'Tourism Management', 'level': ""Bachelor's Degree"", 'cgpa': '3.38', 'university': 'TUNKU ABDUL RAHMAN UNIVERSITY COLLEGE', 'start_date': '2016', 'year_of_graduation': '2016'}, {'field_of_study': 'Hospitality Management', 'level': 'Diploma', 'cgpa': '3.20', 'university': 'N/A', 'start_date': '2014', 'year_of_graduation': '2014'}]"

# Professional certificates
$ws.Range("J2").Value = "['Googe Analytics for Beginners', 'The Fundamental SQL Bootcamp', 'Python Programming for Beginners', 'Mastering SQL server']"

# Skill group
$ws.Range("K2").Value = "['Written and verbal communications in Chinese, English & Malay', 'Familiar with airlines reservation system, Sabre & Amadeus', 'Experienced with hotel distribution channels', 'Quick learner, ability to learn new skills quickly, act on feedback constructively and apply new knowledge immediately with the ability to identify learning opportunities']"

# Language
$ws.Range("L2").Value = "['Chinese', 'English', 'Malay']"
